# The commit moves the second data block (rows 3:4 -> values 5..12) from
# Sheet2 over to Sheet3 (A1:D2), and leaves Sheet2 as the active sheet
# with the now-empty source range selected. This mirrors a user
# selecting Sheet2!A3:D4, cutting it, switching to Sheet3 and pasting
# it at A1.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Cut A3:D4 on Sheet2 and paste it at A1 on Sheet3.
$ws2.Range("A3:D4").Cut($ws3.Range("A1")) | Out-Null

# Leave the pasted range selected on Sheet3 ...
$ws3.Range("A1:D2").Select() | Out-Null

# ... then return to Sheet2, which ends up as the active/selected tab,
# with the vacated source range still selected.
$ws2.Activate() | Out-Null
$ws2.Range("A3:D4").Select() | Out-Null
